$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stimuli order: image (col B), word (col C), category (col D) per row,
# rebuilding the shared-string table that the commit replaced wholesale.
$ws.Range("B2").Value = "dog/dog075.png"
$ws.Range("C2").Value = "bitten"
$ws.Range("D2").Value = "dog"
$ws.Range("B3").Value = "dog/dog105.png"
$ws.Range("C3").Value = "opfern"
$ws.Range("D3").Value = "dog"
$ws.Range("B4").Value = "dog/dog080.png"
$ws.Range("C4").Value = "laufen"
$ws.Range("D4").Value = "dog"
$ws.Range("B5").Value = "dog/dog065.png"
$ws.Range("C5").Value = "kehren"
$ws.Range("D5").Value = "dog"
$ws.Range("B6").Value = "dog/dog111.png"
$ws.Range("C6").Value = "jubeln"
$ws.Range("D6").Value = "dog"
$ws.Range("B7").Value = "car/car065.png"
$ws.Range("C7").Value = "haken"
$ws.Range("D7").Value = "car"
$ws.Range("B8").Value = "dog/dog085.png"
$ws.Range("C8").Value = "formen"
$ws.Range("D8").Value = "dog"
$ws.Range("B9").Value = "dog/dog073.png"
$ws.Range("C9").Value = "strahlen"
$ws.Range("D9").Value = "dog"
$ws.Range("B10").Value = "car/car064.png"
$ws.Range("C10").Value = "runden"
$ws.Range("D10").Value = "car"
$ws.Range("B11").Value = "dog/dog093.png"
$ws.Range("C11").Value = "backen"
$ws.Range("D11").Value = "dog"
$ws.Range("B12").Value = "car/car091.png"
$ws.Range("C12").Value = "fliehen"
$ws.Range("D12").Value = "car"
$ws.Range("B13").Value = "car/car106.png"
$ws.Range("C13").Value = "kaufen"
$ws.Range("D13").Value = "car"
$ws.Range("B14").Value = "car/car107.png"
$ws.Range("C14").Value = "schmecken"
$ws.Range("D14").Value = "car"
$ws.Range("B15").Value = "dog/dog074.png"
$ws.Range("C15").Value = "ehren"
$ws.Range("D15").Value = "dog"
$ws.Range("B16").Value = "dog/dog064.png"
$ws.Range("C16").Value = "tagen"
$ws.Range("D16").Value = "dog"
$ws.Range("B17").Value = "dog/dog107.png"
$ws.Range("C17").Value = "raten"
$ws.Range("D17").Value = "dog"
$ws.Range("B18").Value = "car/car082.png"
$ws.Range("C18").Value = "liefern"
$ws.Range("D18").Value = "car"
$ws.Range("B19").Value = "car/car088.png"
$ws.Range("C19").Value = "währen"
$ws.Range("D19").Value = "car"
$ws.Range("B20").Value = "dog/dog113.png"
$ws.Range("C20").Value = "bleiben"
$ws.Range("D20").Value = "dog"
$ws.Range("B21").Value = "car/car087.png"
$ws.Range("C21").Value = "enden"
$ws.Range("D21").Value = "car"
$ws.Range("B22").Value = "dog/dog090.png"
$ws.Range("C22").Value = "schätzen"
$ws.Range("D22").Value = "dog"
$ws.Range("B23").Value = "car/car094.png"
$ws.Range("C23").Value = "füttern"
$ws.Range("D23").Value = "car"
$ws.Range("B24").Value = "car/car120.png"
$ws.Range("C24").Value = "hoffen"
$ws.Range("D24").Value = "car"
$ws.Range("B25").Value = "car/car096.png"
$ws.Range("C25").Value = "spielen"
$ws.Range("D25").Value = "car"
$ws.Range("B26").Value = "dog/dog066.png"
$ws.Range("C26").Value = "gelten"
$ws.Range("D26").Value = "dog"
$ws.Range("B27").Value = "dog/dog099.png"
$ws.Range("C27").Value = "wiegen"
$ws.Range("D27").Value = "dog"
$ws.Range("B28").Value = "car/car086.png"
$ws.Range("C28").Value = "drohen"
$ws.Range("D28").Value = "car"
$ws.Range("B29").Value = "dog/dog120.png"
$ws.Range("C29").Value = "posten"
$ws.Range("D29").Value = "dog"
$ws.Range("B30").Value = "car/car071.png"
$ws.Range("C30").Value = "sieben"
$ws.Range("D30").Value = "car"
$ws.Range("B31").Value = "car/car092.png"
$ws.Range("C31").Value = "wenden"
$ws.Range("D31").Value = "car"
$ws.Range("B32").Value = "car/car113.png"
$ws.Range("C32").Value = "saufen"
$ws.Range("D32").Value = "car"
$ws.Range("B33").Value = "car/car110.png"
$ws.Range("C33").Value = "nehmen"
$ws.Range("D33").Value = "car"
